# "image cropping and size fixes"
#
# The slide has a single content picture ("Picture 4", shape id=5) whose
# bottom edge gets cropped a little (srcRect b="650", i.e. 0.65%) and whose
# displayed height shrinks to match (4000500 EMU -> 3974523 EMU), while its
# width/left/top stay untouched.
#
# PowerPoint's PictureFormat.CropLeft/Top/Right/Bottom properties are
# expressed in points against the *source image* treated at 96 dpi
# (600 x 315 px -> 450 x 236.25 "pt"), independent of the points used for
# the shape's own Left/Top/Width/Height (which map 1 px -> 1 pt for this
# 72 dpi jpeg). Cropping 0.65% off the bottom of a 236.25-"pt"-tall source
# is 1.535625 pt of CropBottom.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 5) {
        $pic = $candidate
    }
}

$pic.PictureFormat.CropBottom = 1.535625
$pic.Height = 312.9546

# --- incidental: the datetimeFigureOut placeholders on the slide master and
# every slide layout get their cached text refreshed to the save date
# (4/30/2020) any time the deck is resaved, regardless of which shape was
# actually edited.
$dateText = "4/30/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDate = $true
                }
            } catch {
                $isDate = $false
            }
            if ($isDate) {
                $shp.TextFrame.TextRange.Text = $dateText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
